# Auto-generated edit script: updates the symbol list per the commit diff
# ("Updated symbol list on Mon Feb 13 04:32:50 UTC 2023 with GitHub Actions").
#
# Columns B/C hold plain text (coin name / link) and are assigned directly.
# Columns D/E hold numeric-looking strings (prices, percentages). The source
# workbook stores every data cell as literal text, so each such value is
# prefixed with a leading apostrophe -- Excel's standard 'force text' marker --
# to stop Excel from auto-converting it into a Number cell on assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "'312.87"),
    @("E2", "'1.63%"),
    @("D3", "'40.08"),
    @("E3", "'-2.44%"),
    @("D4", "'5.211"),
    @("E4", "'-0.46%"),
    @("D5", "'0.07585"),
    @("E5", "'-1.02%"),
    @("B6", "GateToken"),
    @("C6", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"),
    @("D6", "'4.332"),
    @("E6", "'0.31%"),
    @("B7", "FTXToken"),
    @("C7", "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"),
    @("D7", "'1.680"),
    @("E7", "'2.68%"),
    @("B8", "MXToken"),
    @("C8", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"),
    @("D8", "'0.9258"),
    @("E8", "'0.97%"),
    @("B9", "BTSEToken"),
    @("C9", "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"),
    @("D9", "'2.424"),
    @("E9", "'-0.21%"),
    @("B10", "LiechtensteinCryptoassetsExchange"),
    @("C10", "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"),
    @("D10", "'0.1203"),
    @("E10", "'-4.14%"),
    @("B11", "WazirX"),
    @("C11", "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"),
    @("D11", "'0.1800"),
    @("E11", "'-2.05%"),
    @("B12", "MandalaExchangeToken"),
    @("C12", "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"),
    @("D12", "'0.09052"),
    @("E12", "'-1.23%"),
    @("B13", "BitrueCoin"),
    @("C13", "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"),
    @("D13", "'0.04133"),
    @("E13", "'-2.93%"),
    @("B14", "BitMartToken"),
    @("C14", "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"),
    @("D14", "'0.1052"),
    @("E14", "'-0.01%"),
    @("B15", "BitForexToken"),
    @("C15", "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"),
    @("D15", "'0.001280"),
    @("E15", "'1.39%"),
    @("B16", "TigerCash"),
    @("C16", "https://coinranking.com/coin/6hIn06L2+tigercash-tch"),
    @("D16", "'0.005830"),
    @("E16", "'-0.89%"),
    @("B17", "UpBots"),
    @("C17", "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"),
    @("D17", "'0.007522"),
    @("E17", "'0.18%"),
    @("B18", "LEO"),
    @("C18", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"),
    @("D18", "'3.349"),
    @("E18", "'-0.21%"),
    @("E19", "'0.55%"),
    @("D20", "'7.638"),
    @("E20", "'5.68%"),
    @("D21", "'0.1349"),
    @("E21", "'-2.62%"),
    @("D22", "'0.2867"),
    @("E22", "'-1.05%"),
    @("D23", "'0.04010"),
    @("E23", "'-1.66%"),
    @("D24", "'0.001269"),
    @("E24", "'0.69%"),
    @("D25", "'0.004052"),
    @("E25", "'-2.30%"),
    @("D26", "'0.0001271"),
    @("D38", "'0.02437"),
    @("E38", "'-1.30%"),
    @("D39", "'0.05174"),
    @("E39", "'-2.08%"),
    @("D40", "'0.007753"),
    @("E40", "'-1.29%"),
    @("D41", "'0.1296"),
    @("E41", "'-1.34%"),
    @("D42", "'0.007605"),
    @("E42", "'11.37%"),
    @("D43", "'0.002897"),
    @("E43", "'51.62%"),
    @("D44", "'0.008205"),
    @("E44", "'5.34%"),
    @("D45", "'0.3122"),
    @("E45", "'2.20%"),
    @("D46", "'0.00006574"),
    @("E46", "'-2.24%"),
    @("D47", "'0.00000000749"),
    @("E47", "'-0.41%"),
    @("D48", "'0.2692"),
    @("E48", "'58.49%"),
    @("D49", "'0.004195"),
    @("E49", "'2.32%"),
    @("D50", "'0.00002097"),
    @("E50", "'-0.41%"),
    @("D51", "'0.0001998"),
    @("E51", "'-0.41%")
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
